$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header cell G4: "assignment" -> "assignment`n(width:lsb:sequence:step)", wrapped ---
$ws.Range("G4").Value = "assignment`n(width:lsb:sequence:step)"
$ws.Range("G4").WrapText = $true

# --- 2. mac_test_array comment (L21): replace "width:lsb:sequence:step" with the long description ---
$ws.Range("L21").Value = "每一路I2S发送的数据包的总和，格式是：width:lsb:sequence:step，step=32相邻两个counter的起始位间隔32b，sequence=16表示16路I2S，lsb=0表示生成的寄存器从bit0开始算起，width=32表示每个counter位宽32b"

# --- 3. Extend the "mac" row-label merge from A16:A20 down to A16:A21 ---
#     (also clears the now-redundant A21 cell, so row 21's span becomes "2:12")
$ws.Range("A16:A21").Merge()

# --- 4. Give the "offset address" column (C) on the new rows the same left-aligned
#        formatting used throughout the rest of the table (matches style used by C5:C21) ---
$ws.Range("C22:C25").HorizontalAlignment = -4131

# --- 5. Add the new "i2s_out" register block as rows 22-25 ---

# Row 22: i2s_out_tdm_num
$ws.Cells.Item(22, 2).Value = "i2s_out_tdm_num"
$ws.Cells.Item(22, 3).Value = "0x1000"
$ws.Cells.Item(22, 7).Value = "32:0:16:256"
$ws.Cells.Item(22, 8).Value = "rw"
$ws.Cells.Item(22, 9).Value = 2
$ws.Cells.Item(22, 12).Value = "tdm值，合法值为2/4/8/16"

# Row 23: i2s_out_is_master
$ws.Cells.Item(23, 2).Value = "i2s_out_is_master"
$ws.Cells.Item(23, 3).Value = "0x1004"
$ws.Cells.Item(23, 7).Value = "32:0:16:256"
$ws.Cells.Item(23, 8).Value = "rw"
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 12).Value = "是否master模式，如是，由FPGA提供时钟信号，寄存器取值范围:0/1"

# Row 24: i2s_out_enable
$ws.Cells.Item(24, 2).Value = "i2s_out_enable"
$ws.Cells.Item(24, 3).Value = "0x1008"
$ws.Cells.Item(24, 7).Value = "32:0:16:256"
$ws.Cells.Item(24, 8).Value = "rw"
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 12).Value = "是否使能，1：使能，0：不使能"

# Row 25: i2s_out_dst fpga index
$ws.Cells.Item(25, 2).Value = "i2s_out_dst fpga index"
$ws.Cells.Item(25, 3).Value = "0x100C"
$ws.Cells.Item(25, 7).Value = "32:0:16:256"
$ws.Cells.Item(25, 8).Value = "rw"
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 12).Value = "目的FPGA的index，1~8，由板子上的拨码开关确定"
